$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 2.67
$ws.Range("I4").Value = 2.57
$ws.Range("M4").Value = 1.01
$ws.Range("N4").Value = 8
$ws.Range("O4").Value = 1.31
$ws.Range("P4").Value = 2.9
$ws.Range("Q4").Value = 1.95
$ws.Range("R4").Value = 1.75
$ws.Range("U4").Value = 1.7
$ws.Range("V4").Value = 1.93
$ws.Range("W4").Value = 8.25
$ws.Range("X4").Value = 13.5
$ws.Range("AB4").Value = 32
$ws.Range("AE4").Value = 13
$ws.Range("AF4").Value = 60
$ws.Range("AG4").Value = 450
$ws.Range("AJ4").Value = 9.5
$ws.Range("AM4").Value = 29
$ws.Range("AN4").Value = 4.65
$ws.Range("AP4").Value = 20
$ws.Range("AQ4").Value = 65
$ws.Range("AS4").Value = 250
$ws.Range("AT4").Value = 2.6
$ws.Range("AU4").Value = 6.4
$ws.Range("AX4").Value = 13.5
$ws.Range("AY4").Value = 19
$ws.Range("BA4").Value = 80
